$wb = $excel.ActiveWorkbook

# Add a new worksheet "postLogin" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "postLogin"

# Header row
$ws.Range("A1").Value = "user_id"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "Expected Code"
$ws.Range("D1").Value = "Description"

# Row 2
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "SacredGroves@FT!@#007"
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = "Correct"

# Row 3
$ws.Range("B3").Value = "SacredGroves@FT!@#007"
$ws.Range("C3").Value = 500
$ws.Range("D3").Value = "empty string userid"

# Row 4
$ws.Range("A4").Value = "Admin"
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = "empty string password"

# Row 5
$ws.Range("A5").Value = "Admin"
$ws.Range("B5").Value = 'Incorrect123 @#$'
$ws.Range("C5").Value = 500
$ws.Range("D5").Value = "Incorrect password"

# Row 6
$ws.Range("A6").Value = "Adm in"
$ws.Range("B6").Value = "SacredGroves@FT!@#007"
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = "space in user_id"

# Row 7
$ws.Range("A7").Value = "Ad*min"
$ws.Range("B7").Value = "SacredGroves@FT!@#007"
$ws.Range("C7").Value = 500
$ws.Range("D7").Value = "Special chars in userid"

# Row 8
$ws.Range("A8").Value = '@#$%^&*()><?.'',;'
$ws.Range("B8").Value = "SacredGroves@FT!@#007"
$ws.Range("C8").Value = 500
$ws.Range("D8").Value = "Special chars in userid"

# Row 9
$ws.Range("B9").Value = "SacredGroves@FT!@#007"
$ws.Range("C9").Value = 400
$ws.Range("D9").Value = "blank user id"

# Row 10
$ws.Range("A10").Value = "Admin"
$ws.Range("C10").Value = 400
$ws.Range("D10").Value = "blank password"

# Match the final selection state shown in the target worksheet
$ws.Range("H11").Select()
